$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @("ECs",  "Hspg2", "Itga2", "ECs",  3, 1, 117.9639543333333, 353.891863,        0.2661690114309019,  0.2661690114309019,  2, 0.6666666666666666, 1.888921,          5.666763, 0.4551532417350329, 0.4551532417350328, 222.8245905832743, 2005.421315249469, 0.121147688402184,   0.121147688402184),
    @("ECs",  "Hspg2", "Itga2", "FAPs", 3, 1, 117.9639543333333, 353.891863,        0.2661690114309019,  0.2661690114309019,  3, 1,                 1.087098333333333, 3.261295, 0.2619465454094788, 0.2619465454094787, 128.2384181491761, 1154.145763342585, 0.06972205303938082, 0.06972205303938081),
    @("ECs",  "Hspg2", "Itga2", "sCs",  3, 1, 117.9639543333333, 353.891863,        0.2661690114309019,  0.2661690114309019,  3, 1,                 1.174057666666667, 3.522173, 0.2829002128554884, 0.2829002128554884, 138.4964849753665, 1246.468364778299, 0.07529926998933707, 0.07529926998933707),
    @("FAPs", "Hspg2", "Itga2", "ECs",  3, 1, 282.6413673333333, 847.9241019999999, 0.6377403483780447,  0.6377403483780446,  2, 0.6666666666666666, 1.888921,          5.666763, 0.4551532417350329, 0.4551532417350328, 533.8872142246473, 4804.984928021825, 0.2902695869494963,  0.2902695869494962),
    @("FAPs", "Hspg2", "Itga2", "FAPs", 3, 1, 282.6413673333333, 847.9241019999999, 0.6377403483780447,  0.6377403483780446,  3, 1,                 1.087098333333333, 3.261295, 0.2619465454094788, 0.2619465454094787, 307.2589593591211, 2765.330634232089, 0.1670538811258663,  0.1670538811258662),
    @("FAPs", "Hspg2", "Itga2", "sCs",  3, 1, 282.6413673333333, 847.9241019999999, 0.6377403483780447,  0.6377403483780446,  3, 1,                 1.174057666666667, 3.522173, 0.2829002128554884, 0.2829002128554884, 331.8372642348496, 2986.535378113646, 0.1804168803026822,  0.1804168803026821),
    @("sCs",  "Hspg2", "Itga2", "ECs",  3, 1, 42.586595,         127.759785,        0.09609064019105341, 0.09609064019105343, 2, 0.6666666666666666, 1.888921,          5.666763, 0.4551532417350329, 0.4551532417350328, 80.44271361399498, 723.9844225259549, 0.0437359663833526,  0.0437359663833526),
    @("sCs",  "Hspg2", "Itga2", "FAPs", 3, 1, 42.586595,         127.759785,        0.09609064019105341, 0.09609064019105343, 3, 1,                 1.087098333333333, 3.261295, 0.2619465454094788, 0.2619465454094787, 46.29581644684166, 416.6623480215749, 0.02517061124423166, 0.02517061124423165),
    @("sCs",  "Hspg2", "Itga2", "sCs",  3, 1, 42.586595,         127.759785,        0.09609064019105341, 0.09609064019105343, 3, 1,                 1.174057666666667, 3.522173, 0.2829002128554884, 0.2829002128554884, 49.99911835697833, 449.992065212805, 0.02718406256346916, 0.02718406256346916)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $values[$c]
    }
}
